# Auto-generated script to update market-price columns (H-N) across multiple sheets
# per the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 588.8333
$ws.Range("I9").Value = 584.625
$ws.Range("K9").Value = 584.625
$ws.Range("M9").Value = -415.625
$ws.Range("H33").Value = 2786.3215
$ws.Range("I33").Value = 2909.739
$ws.Range("K33").Value = 2909.739
$ws.Range("M33").Value = -2680.739
$ws.Range("H43").Value = 162665.61
$ws.Range("J43").Value = 209060.4
$ws.Range("L43").Value = 209060.4
$ws.Range("N43").Value = -209198.4
$ws.Range("H63").Value = 85500
$ws.Range("J63").Value = 116000
$ws.Range("L63").Value = 116000
$ws.Range("N63").Value = -117248
$ws.Range("H66").Value = 85500
$ws.Range("J66").Value = 116000
$ws.Range("L66").Value = 348000
$ws.Range("N66").Value = -354240
$ws.Range("H69").Value = 44169
$ws.Range("I69").Value = 26666.334
$ws.Range("K69").Value = 79999.00199999999
$ws.Range("M69").Value = -79125.00199999999
$ws.Range("H72").Value = 44169
$ws.Range("I72").Value = 26666.334
$ws.Range("K72").Value = 239997.006
$ws.Range("M72").Value = -235629.006
$ws.Range("H101").Value = 1344
$ws.Range("I101").Value = 1266.125
$ws.Range("K101").Value = 3798.375
$ws.Range("M101").Value = -2176.375
$ws.Range("H116").Value = 6533
$ws.Range("I116").Value = 6429.636
$ws.Range("J116").Value = 6675.125
$ws.Range("K116").Value = 6429.636
$ws.Range("L116").Value = 6675.125
$ws.Range("M116").Value = -2987.636
$ws.Range("N116").Value = -13559.125
$ws.Range("H123").Value = 69990
$ws.Range("J123").Value = 69990
$ws.Range("L123").Value = 69990
$ws.Range("N123").Value = -79790
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 71429910
$ws.Range("I45").Value = 100001010
$ws.Range("J45").Value = 2165
$ws.Range("K45").Value = 100001010
$ws.Range("L45").Value = 2165
$ws.Range("M45").Value = -100000633
$ws.Range("N45").Value = -2919
$ws.Range("H74").Value = 7227370
$ws.Range("I74").Value = 8931086
$ws.Range("J74").Value = 1264362.8
$ws.Range("K74").Value = 8931086
$ws.Range("L74").Value = 1264362.8
$ws.Range("M74").Value = -8930212
$ws.Range("N74").Value = -1266110.8
$ws.Range("H77").Value = 7227370
$ws.Range("I77").Value = 8931086
$ws.Range("J77").Value = 1264362.8
$ws.Range("K77").Value = 44655430
$ws.Range("L77").Value = 6321814
$ws.Range("M77").Value = -44651062
$ws.Range("N77").Value = -6330550
$ws.Range("H81").Value = 92329.664
$ws.Range("J81").Value = 92329.664
$ws.Range("L81").Value = 92329.664
$ws.Range("N81").Value = -94325.664
$ws.Range("H84").Value = 92329.664
$ws.Range("J84").Value = 92329.664
$ws.Range("L84").Value = 276988.992
$ws.Range("N84").Value = -286972.992
$ws.Range("H95").Value = 62494.75
$ws.Range("J95").Value = 62494.75
$ws.Range("L95").Value = 62494.75
$ws.Range("N95").Value = -67986.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 37500
$ws.Range("I54").Value = 37500
$ws.Range("K54").Value = 37500
$ws.Range("M54").Value = -37016

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 408966.22
$ws.Range("I31").Value = 8869.147999999999
$ws.Range("K31").Value = 8869.147999999999
$ws.Range("M31").Value = -8574.147999999999
$ws.Range("H34").Value = 408966.22
$ws.Range("I34").Value = 8869.147999999999
$ws.Range("K34").Value = 8869.147999999999
$ws.Range("M34").Value = -8667.147999999999
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 30000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -31316
$ws.Range("H88").Value = 28484.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 28484.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 28484.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -29296.5
$ws.Range("H91").Value = 28484.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 28484.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 28484.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -31292.5
$ws.Range("H121").Value = 29399.8
$ws.Range("I121").Value = 19999
$ws.Range("J121").Value = 31750
$ws.Range("K121").Value = 19999
$ws.Range("L121").Value = 31750
$ws.Range("M121").Value = -18689
$ws.Range("N121").Value = -34370

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3899.923
$ws.Range("I88").Value = 2900
$ws.Range("J88").Value = 4081.7273
$ws.Range("K88").Value = 8700
$ws.Range("L88").Value = 12245.1819
$ws.Range("M88").Value = -8272
$ws.Range("N88").Value = -13101.1819
$ws.Range("H91").Value = 3899.923
$ws.Range("I91").Value = 2900
$ws.Range("J91").Value = 4081.7273
$ws.Range("K91").Value = 8700
$ws.Range("L91").Value = 12245.1819
$ws.Range("M91").Value = -7218
$ws.Range("N91").Value = -15209.1819
$ws.Range("H103").Value = 648.2222
$ws.Range("I103").Value = 224
$ws.Range("K103").Value = 672
$ws.Range("M103").Value = 207
$ws.Range("H113").Value = 2662
$ws.Range("I113").Value = 876.5
$ws.Range("J113").Value = 3257.1667
$ws.Range("K113").Value = 2629.5
$ws.Range("L113").Value = 9771.500100000001
$ws.Range("M113").Value = -459.5
$ws.Range("N113").Value = -14111.5001
$ws.Range("H131").Value = 4313.097
$ws.Range("J131").Value = 3851.541
$ws.Range("L131").Value = 11554.623
$ws.Range("N131").Value = -21634.623

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 33147.5
$ws.Range("J24").Value = 33147.5
$ws.Range("L24").Value = 33147.5
$ws.Range("N24").Value = -33493.5
$ws.Range("H80").Value = 2895.7144
$ws.Range("I80").Value = 2705.2
$ws.Range("K80").Value = 2705.2
$ws.Range("M80").Value = -1707.2
$ws.Range("H83").Value = 2895.7144
$ws.Range("I83").Value = 2705.2
$ws.Range("K83").Value = 13526
$ws.Range("M83").Value = -8534
$ws.Range("H93").Value = 59790
$ws.Range("J93").Value = 59790
$ws.Range("L93").Value = 59790
$ws.Range("N93").Value = -63534
$ws.Range("H97").Value = 1307.9375
$ws.Range("I97").Value = 1228.4667
$ws.Range("K97").Value = 1228.4667
$ws.Range("M97").Value = -732.4666999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 37495
$ws.Range("J48").Value = 37495
$ws.Range("L48").Value = 37495
$ws.Range("N48").Value = -38817
$ws.Range("H68").Value = 2746.8333
$ws.Range("J68").Value = 2998.25
$ws.Range("L68").Value = 2998.25
$ws.Range("N68").Value = -4496.25
$ws.Range("H71").Value = 2746.8333
$ws.Range("J71").Value = 2998.25
$ws.Range("L71").Value = 14991.25
$ws.Range("N71").Value = -22479.25
$ws.Range("H82").Value = 1136.25
$ws.Range("I82").Value = 1063.7858
$ws.Range("J82").Value = 1305.3334
$ws.Range("K82").Value = 1063.7858
$ws.Range("L82").Value = 1305.3334
$ws.Range("M82").Value = -702.7858000000001
$ws.Range("N82").Value = -2027.3334
$ws.Range("H85").Value = 1136.25
$ws.Range("I85").Value = 1063.7858
$ws.Range("J85").Value = 1305.3334
$ws.Range("K85").Value = 1063.7858
$ws.Range("L85").Value = 1305.3334
$ws.Range("M85").Value = 184.2141999999999
$ws.Range("N85").Value = -3801.3334
$ws.Range("H93").Value = 52632930
$ws.Range("I93").Value = 66668084
$ws.Range("K93").Value = 66668084
$ws.Range("M93").Value = -66666836
